$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (rows 2-43) from 45727 to 45728 (i.e. +1 day)
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45727) {
        $cell.Value2 = 45728
    }
}
